$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.357.78"
$ws.Range("D3").Value = "3.686.15"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "680.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.438"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "4.308.02"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "3.708.68"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "69.321.90"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "80.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "3.834.83"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("D35").Value = "3.675.01"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.162"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "169.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000277"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("E51").Value = "  -0.34%  "
